$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.003.50'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.301.54'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.45%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.42'
$ws.Range("E5").Value = '  -2.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.50'
$ws.Range("E6").Value = '  +1.66%  '
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.607'
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.72'
$ws.Range("E10").Value = '  -0.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0906'
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.26'
$ws.Range("E12").Value = '  -2.64%  '
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.999'
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.28'
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.648.26'
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.298.62'
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.828.32'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.29'
$ws.Range("E19").Value = '  -4.19%  '
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.34'
$ws.Range("E21").Value = '  -0.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.34'
$ws.Range("E22").Value = '  -0.83%  '
$ws.Range("E23").Value = '  -2.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.08'
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.24'
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.97'
$ws.Range("E27").Value = '  +0.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.25'
$ws.Range("E28").Value = '  +10.23%  '
$ws.Range("E29").Value = '  -2.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.23'
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.22'
$ws.Range("E31").Value = '  -2.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '164.96'
$ws.Range("E32").Value = '  -0.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0859'
$ws.Range("E33").Value = '  -2.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.84'
$ws.Range("E34").Value = '  +5.05%  '
$ws.Range("E35").Value = '  -0.67%  '
$ws.Range("E36").Value = '  -1.88%  '
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0349'
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("E39").Value = '  +3.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.63'
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '107.54'
$ws.Range("E41").Value = '  +10.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.58'
$ws.Range("E42").Value = '  -2.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '71.49'
$ws.Range("E43").Value = '  +2.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.228'
$ws.Range("E44").Value = '  +1.40%  '
$ws.Range("E45").Value = '  +0.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.33'
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.721.87'
$ws.Range("E47").Value = '  +4.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '111.75'
$ws.Range("E48").Value = '  -4.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.67'
$ws.Range("E49").Value = '  -5.12%  '
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("E51").Value = '  -2.47%  '
